$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 122; this shifts the existing rows
# 122-222 down to 123-223 (and extends dimension to R223).
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new record.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg,
# F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
# K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg,
# Q Kg o Unidades, R Clasificacion.
$ws.Cells.Item(122, 1).Value = 9
$ws.Cells.Item(122, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(122, 3).Value = "Metropolitana"
$ws.Cells.Item(122, 4).NumberFormat = $ws.Cells.Item(123, 4).NumberFormat
$ws.Cells.Item(122, 4).Value = 44574
$ws.Cells.Item(122, 5).Value = 13
$ws.Cells.Item(122, 6).Value = 100112021
$ws.Cells.Item(122, 7).Value = "Ají"
$ws.Cells.Item(122, 8).Value = "Americana (o)"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 34
$ws.Cells.Item(122, 11).Value = 20000
$ws.Cells.Item(122, 12).Value = 22000
$ws.Cells.Item(122, 13).Value = 21000
$ws.Cells.Item(122, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(122, 15).Value = "Región Metropolitana"
$ws.Cells.Item(122, 16).Value = 840
$ws.Cells.Item(122, 17).Value = 25
$ws.Cells.Item(122, 18).Value = "Hortaliza"
